$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 190, shifting existing rows 190:248 down to 191:249
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row 190 with the new weekly data record
$ws.Range("A190").Value = 4
$ws.Range("B190").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C190").Value = "Los Lagos"
$ws.Range("D190").Value = 44809
$ws.Range("E190").Value = 10
$ws.Range("F190").Value = "Fruta"
$ws.Range("G190").Value = 100101
$ws.Range("H190").Value = "Berries"
$ws.Range("I190").Value = 100112025
$ws.Range("J190").Value = "Frutilla"
$ws.Range("K190").Value = "Sin especificar"
$ws.Range("L190").Value = "Segunda"
$ws.Range("M190").Value = 200
$ws.Range("N190").Value = 15000
$ws.Range("O190").Value = 16000
$ws.Range("P190").Value = 15500
$ws.Range("Q190").Value = "$/bandeja 7 kilos"
$ws.Range("R190").Value = "Provincia de Melipilla"
$ws.Range("S190").Value = 2214
$ws.Range("T190").Value = 7
